# edit.ps1 - PowerShell COM-interop edit script
#
# Reproduces the target diff:
#   1) Re-style the three tables (on slides 14, 15, 16) from table style
#      {F8B178E9-22AB-4F3B-A4C9-92F57F457300} to
#      {286298FD-098D-45F2-AF25-1EB7391DAF29}.
#   2) Swap the two presentation themes: the slide-master theme
#      (ppt/theme/theme1.xml, "Integral"/"Red Violet") takes on the
#      colours that used to live in the notes-master theme
#      (ppt/theme/theme2.xml, "Office Theme").

$p = $ppt.ActivePresentation

# --- 1) Update table style IDs on the three slides that contain tables ---
$newStyleId = "{286298FD-098D-45F2-AF25-1EB7391DAF29}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId, $true)
        }
    }
}

# --- 2) Re-colour the (shared) presentation theme to the "Office Theme" ---
#        palette that previously lived in the secondary theme part.
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

# Index -> RGB (Windows COLORREF / BGR-packed long), in clrScheme order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @{
    1  = 0x000000   # dk1      000000
    2  = 0xFFFFFF   # lt1      FFFFFF
    3  = 0x6A5444   # dk2      44546A (BGR-packed)
    4  = 0xE6E6E7   # lt2      E7E6E6 (BGR-packed)
    5  = 0xD59B5B   # accent1  5B9BD5 (BGR-packed)
    6  = 0x317DED   # accent2  ED7D31 (BGR-packed)
    7  = 0xA5A5A5   # accent3  A5A5A5
    8  = 0x00C0FF   # accent4  FFC000 (BGR-packed)
    9  = 0xC47244   # accent5  4472C4 (BGR-packed)
    10 = 0x47AD70   # accent6  70AD47 (BGR-packed)
    11 = 0xC16305   # hlink    0563C1 (BGR-packed)
    12 = 0x724F95   # folHlink 954F72 (BGR-packed)
}

foreach ($idx in 1..12) {
    $colorScheme.Colors($idx).RGB = $officeThemeColors[$idx]
}
